$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A updates: form-field id naming scheme changed from underscore to
# dash separated prefixes (a1_/o1_/n1_/z1_), and system_submitter_id was
# renamed to SubmitterID (RCT / LS encrypted submitter id dictionary entry).
$ws.Range("A5").Value = "SubmitterID"
$ws.Range("A6").Value = "a1-pid"
$ws.Range("A7").Value = "a1-fid"
$ws.Range("A8").Value = "a1-enrolfacility"
$ws.Range("A9").Value = "a1-enroldate"
$ws.Range("A10").Value = "a1-relationship"
$ws.Range("A11").Value = "a1-name"
$ws.Range("A12").Value = "a1-caregiver"
$ws.Range("A13").Value = "a1-phonenb"
$ws.Range("A14").Value = "a1-contact-success"
$ws.Range("A15").Value = "a1-contact-a4_d_1a"
$ws.Range("A16").Value = "a1-contact-a4_d_1b"
$ws.Range("A18").Value = "o1-o1_1a"
$ws.Range("A19").Value = "o1-o1_2"
$ws.Range("A20").Value = "o1-o1_1"
$ws.Range("A21").Value = "n1-o3_1"
$ws.Range("A22").Value = "n1-o3_1a"
$ws.Range("A23").Value = "n1-o3_1a_o"
$ws.Range("A24").Value = "n1-n1_4"
$ws.Range("A25").Value = "n1-n1_3"
$ws.Range("A26").Value = "n1-ref_location_name"
$ws.Range("A27").Value = "n1-n1_3o"
$ws.Range("A28").Value = "n1-n1_2b"
$ws.Range("A29").Value = "n1-n1_2o"
$ws.Range("A30").Value = "n1-n1_5"
$ws.Range("A31").Value = "n1-n1_6"
$ws.Range("A32").Value = "n1-n1_7"
$ws.Range("A33").Value = "n1-maxduration"
$ws.Range("A34").Value = "n1-n1_8a"
$ws.Range("A35").Value = "n1-n1_8"
$ws.Range("A36").Value = "z1-qual"

# Update the view: scroll back to top-left and move the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F9").Select()
